# Updates the cryptos list (prices / 1h volume %) per the Dec 27 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (rows re-ranked) ---
$ws.Range("B30").Value = "WEMIXToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

# --- Volume(1h) percentage text updates (column E) ---
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("E3").Value = "  +6.57%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +8.81%  "
$ws.Range("E6").Value = "  -5.10%  "
$ws.Range("E7").Value = "  +2.68%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("E14").Value = "  +10.47%  "
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("E16").Value = "  +6.78%  "
$ws.Range("E17").Value = "  +9.21%  "
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("E23").Value = "  +15.95%  "
$ws.Range("E24").Value = "  +5.48%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +9.06%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  -3.25%  "
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("E36").Value = "  +3.77%  "
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  +16.59%  "
$ws.Range("E41").Value = "  +13.24%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  +3.84%  "
$ws.Range("E47").Value = "  +8.98%  "
$ws.Range("E48").Value = "  +8.69%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("E51").Value = "  +6.62%  "

# --- Price text updates (column D). These look like numbers to Excel, so
# a leading apostrophe forces literal text, preserving exact digits/
# trailing zeros; Style="Normal" then drops the quote-prefix formatting
# so the cell keeps its original (unstyled) appearance. ---
$ws.Range("D2").Value = "'43.539.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.378.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'323.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'107.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Value = "'0.635"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'42.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "'8.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'16.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "'2.737.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'2.438.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'43.544.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'7.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'75.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'266.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'9.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'12.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'39.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'23.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'3.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'174.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.0924"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'5.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'4.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.131"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'4.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.0372"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Value = "'2.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Value = "'71.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Value = "'12.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'5.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'112.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'9.30"
$ws.Range("D48").Style = "Normal"
